# Juno: check in to OLPRODLOC.
# Localizes the "Sales report" sheet into Indonesian ("Laporan penjualan"):
#  - renames the worksheet
#  - replaces the header / quarter labels with their Indonesian equivalents,
#    re-applying explicit run formatting (Aptos Narrow, 11pt, black) the way
#    a paste/retype of localized text would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet to its Indonesian title.
$ws.Name = "Laporan penjualan"

# 2) Map of cell -> new (translated) text.
$translations = @{
    "A1" = "Kuartal Tahun"
    "B1" = "Barat Tengah"
    "C1" = "Gunung"
    "D1" = "Timur laut"
    "E1" = "Selatan"
    "F1" = "Tenggara"
    "G1" = "Barat"
    "A2" = "Q1 2022"
    "A3" = "Q2 2022"
    "A4" = "Q3 2022"
    "A5" = "Q4 2022"
    "A6" = "Q1-2023"
    "A7" = "Q2-2023"
    "A8" = "Q3-2023"
    "A9" = "Q4-2023"
}

$order = @("A1","B1","C1","D1","E1","F1","G1","A2","A3","A4","A5","A6","A7","A8","A9")

foreach ($addr in $order) {
    $text = $translations[$addr]
    $cell = $ws.Range($addr)
    $cell.Value = $text

    $len = $cell.Value2.Length

    # Re-apply the run formatting in two adjacent chunks (rather than one
    # call spanning the whole string) so the shared-string keeps an explicit
    # rich-text run - matching how the localized strings were authored.
    $firstLen = $len - 1
    if ($firstLen -lt 1) { $firstLen = 1 }

    $head = $cell.Characters(1, $firstLen)
    $head.Font.Name = "Aptos Narrow"
    $head.Font.Size = 11
    $head.Font.Color = 0

    if ($len -gt $firstLen) {
        $tail = $cell.Characters($firstLen + 1, $len - $firstLen)
        $tail.Font.Name = "Aptos Narrow"
        $tail.Font.Size = 11
        $tail.Font.Color = 0
    }
}

Write-Host "Localized sheet to Indonesian: $($ws.Name)"
